$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting rows 30:65 down to 31:66
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with data
$ws.Cells.Item(30, 1).Value = 7
$ws.Cells.Item(30, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(30, 3).Value = "Ñuble"
$ws.Cells.Item(30, 4).Value = 44589
$ws.Cells.Item(30, 5).Value = 16
$ws.Cells.Item(30, 6).Value = 100112031
$ws.Cells.Item(30, 7).Value = "Poroto verde"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 60
$ws.Cells.Item(30, 11).Value = 24000
$ws.Cells.Item(30, 12).Value = 25000
$ws.Cells.Item(30, 13).Value = 24500
$ws.Cells.Item(30, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value = "Región del Maule"
$ws.Cells.Item(30, 16).Value = 980
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
